$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = "13-10-2025"
$newPrice = "The price of gold in India today is ₹12,540 per gram for 24 karat gold, ₹11,495 per gram for 22 karat gold and ₹9,405 per gram for 18 karat gold (also called 999 gold)."

# Copy the previous data row's formatting down to the new row, then set values
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A22").Value = $newDate
$ws.Range("B22").Value = $newPrice

$wb.Save()
